$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "('Angel', ['Token Creature — Angel', 'Flying', '4/4'])"
$ws.Range("A3").Value = "('Beast', ['Token Creature — Beast', '3/3'])"
$ws.Range("A4").Value = "('Myr', ['Token Artifact Creature — Myr', '1/1'])"
$ws.Range("A5").Value = "('Pentavite', ['Token Artifact Creature — Pentavite', 'Flying', '1/1'])"
$ws.Range("A6").Value = "('Powder Keg', ['{2}', 'Artifact', 'At the beginning of your upkeep, you may put a fuse counter on Powder Keg.', '{T}, Sacrifice Powder Keg: Destroy each artifact and creature with converted mana cost equal to the number of fuse counters on Powder Keg.'])"
$ws.Range("A7").Value = "('Spirit', ['Token Creature — Spirit', '1/1'])"

$ws.Range("A8:A23").EntireRow.Delete()
